$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.084.25"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "1.781.64"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.38"
$ws.Range("E5").Value = "  -0.71%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.84"
$ws.Range("E8").Value = "  -1.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0947"
$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").Value = "2.037.83"
$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("D13").Value = "1.801.25"
$ws.Range("E13").Value = "  +0.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.93"
$ws.Range("E14").Value = "  -4.59%  "

$ws.Range("D15").Value = "34.089.29"
$ws.Range("E15").Value = "  +0.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.622"
$ws.Range("E16").Value = "  -0.43%  "

$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.56"
$ws.Range("E18").Value = "  -0.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.71"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("E20").Value = "  +1.59%  "

$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.87"
$ws.Range("E22").Value = "  +0.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.10"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -1.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.67"

$ws.Range("E26").Value = "  -0.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.25"
$ws.Range("E27").Value = "  +0.18%  "

$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("E30").Value = "  -1.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0517"
$ws.Range("E31").Value = "  -0.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.71"
$ws.Range("E32").Value = "  +1.60%  "

$ws.Range("E33").Value = "  +2.00%  "

$ws.Range("E34").Value = "  -2.52%  "

$ws.Range("D35").Value = "1.450.89"
$ws.Range("E35").Value = "  +3.39%  "

$ws.Range("E36").Value = "  +4.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.650"
$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("E38").Value = "  +0.79%  "

$ws.Range("E39").Value = "  -0.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.40"
$ws.Range("E40").Value = "  +1.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "80.67"

$ws.Range("E42").Value = "  +1.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.79"
$ws.Range("E43").Value = "  +0.48%  "

$ws.Range("E44").Value = "  -0.61%  "

$ws.Range("E45").Value = "  +1.79%  "

$ws.Range("E46").Value = "  -1.42%  "

$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0133"
$ws.Range("E48").Value = "  -4.82%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.937.71"
$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.41"
$ws.Range("E50").Value = "  -2.95%  "

$ws.Range("E51").Value = "  +0.20%  "
